$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-snapshot column is inserted right before the "nom" column
# (old AE), pushing "nom" from AE->AF and "url_produit" from AF->AG.
$ws.Range("AE1").EntireColumn.Insert()

# Header for the freshly inserted column: the new snapshot timestamp.
$ws.Range("AE1").Value = "2026-01-29 00:59:01"

# Populate the new snapshot column with the latest known price for each
# product (same value as the previous snapshot column, AD), for every
# data row that actually has a price. Rows without a price (AD blank)
# are left blank, which is already the state after the column insert.
$ws.Range("AE2:AE80").Value = $ws.Range("AD2:AD80").Value()
